$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 7 with the new approver, copying the style from row 6
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A7").Value = "Monica Maria Cardona Suarez (Gerente)"
$ws.Range("B7").Value = 6

$ws.Range("A5").Select()
